# Updates cryptos list values (prices + 1h volume change) to match the
# latest scrape, and fixes the ordering of three coin pairs whose rows
# had swapped (Dai/PancakeSwap, Monero/Toncoin, MultiversX/THORChain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.832.13"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.228.14"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'249.61"
$ws.Range("E5").Value = "  +7.64%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "'71.71"
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +6.84%  "
$ws.Range("D10").Value = "'41.44"
$ws.Range("E10").Value = "  +19.47%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'58.35"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'7.22"
$ws.Range("E13").Value = "  +8.60%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "2.560.24"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "'15.03"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "'0.866"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "2.225.65"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "41.792.05"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'73.21"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'236.20"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  +12.95%  "
$ws.Range("D25").Value = "'4.03"
$ws.Range("E25").Value = "  +10.14%  "
$ws.Range("D28").Value = "'10.75"
$ws.Range("E28").Value = "  +8.31%  "
$ws.Range("D31").Value = "'20.92"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").Value = "'0.125"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'5.59"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").Value = "'0.0731"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'4.73"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "'26.03"
$ws.Range("E37").Value = "  +24.08%  "
$ws.Range("D38").Value = "'3.97"
$ws.Range("E38").Value = "  +10.99%  "
$ws.Range("E39").Value = "  +15.32%  "
$ws.Range("D40").Value = "'2.30"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D43").Value = "'11.97"
$ws.Range("E43").Value = "  +22.51%  "
$ws.Range("D44").Value = "'0.209"
$ws.Range("E44").Value = "  +12.36%  "
$ws.Range("D45").Value = "'4.87"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("D46").Value = "'8.77"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +10.24%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +7.97%  "
$ws.Range("E51").Value = "  +2.52%  "

# Full row replacements (Coin/Link/Price/Volume swapped between pairs)
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.56"
$ws.Range("E27").Value = "  +9.67%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'171.82"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'68.48"
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").Value = "'5.99"
$ws.Range("E42").Value = "  +0.81%  "
